$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing F/G values for rows with revised AgTests/AgPosit figures ---
$ws.Range("F271").Value = 45775
$ws.Range("G271").Value = 1732
$ws.Range("F324").Value = 249805
$ws.Range("F326").Value = 418155
$ws.Range("G326").Value = 3809
$ws.Range("F334").Value = 192788
$ws.Range("F335").Value = 150216
$ws.Range("G335").Value = 3753
$ws.Range("F349").Value = 159211
$ws.Range("G349").Value = 2754
$ws.Range("F350").Value = 126822
$ws.Range("G350").Value = 2780
$ws.Range("F351").Value = 150429
$ws.Range("G351").Value = 2805
$ws.Range("F352").Value = 307121
$ws.Range("G352").Value = 3541
$ws.Range("F355").Value = 221717
$ws.Range("G355").Value = 3432
$ws.Range("F356").Value = 159786
$ws.Range("G356").Value = 2878
$ws.Range("F357").Value = 138021
$ws.Range("G357").Value = 3008
$ws.Range("F358").Value = 158594
$ws.Range("G358").Value = 2603
$ws.Range("F359").Value = 321099
$ws.Range("G359").Value = 3328
$ws.Range("F360").Value = 749621
$ws.Range("G360").Value = 5137
$ws.Range("F362").Value = 228537
$ws.Range("F369").Value = 234695
$ws.Range("G369").Value = 2604
$ws.Range("F372").Value = 178645
$ws.Range("G372").Value = 1854
$ws.Range("F394").Value = 166241
$ws.Range("G394").Value = 634
$ws.Range("F395").Value = 751337
$ws.Range("F399").Value = 201064
$ws.Range("F400").Value = 150805
$ws.Range("G400").Value = 763
$ws.Range("F402").Value = 717010
$ws.Range("G402").Value = 1386
$ws.Range("F403").Value = 351841
$ws.Range("F404").Value = 224660
$ws.Range("G404").Value = 910
$ws.Range("F405").Value = 173845
$ws.Range("G405").Value = 692
$ws.Range("F406").Value = 170803
$ws.Range("G406").Value = 681
$ws.Range("F407").Value = 158108
$ws.Range("G407").Value = 672
$ws.Range("F408").Value = 303816
$ws.Range("G408").Value = 837
$ws.Range("F409").Value = 703103
$ws.Range("F410").Value = 363796
$ws.Range("F413").Value = 148903
$ws.Range("F414").Value = 145643
$ws.Range("F415").Value = 306076
$ws.Range("F416").Value = 651591
$ws.Range("G416").Value = 918
$ws.Range("F417").Value = 329314
$ws.Range("G417").Value = 570
$ws.Range("F418").Value = 199953
$ws.Range("G418").Value = 695
$ws.Range("F419").Value = 146235
$ws.Range("G419").Value = 502
$ws.Range("F420").Value = 135075
$ws.Range("G420").Value = 488
$ws.Range("F421").Value = 144399
$ws.Range("G421").Value = 519

# --- Append new daily rows (2021-04-30, 2021-05-01, 2021-05-02) ---
$ws.Range("A422").Value = 44316
$ws.Range("B422").Value = 382720
$ws.Range("C422").Value = 8431
$ws.Range("D422").Value = 509
$ws.Range("E422").Value = 11732
$ws.Range("F422").Value = 268769
$ws.Range("G422").Value = 574
$ws.Range("A423").Value = 44317
$ws.Range("B423").Value = 383098
$ws.Range("C423").Value = 34967
$ws.Range("D423").Value = 378
$ws.Range("E423").Value = 11766
$ws.Range("F423").Value = 368311
$ws.Range("G423").Value = 585
$ws.Range("A424").Value = 44318
$ws.Range("B424").Value = 383228
$ws.Range("C424").Value = 12176
$ws.Range("D424").Value = 130
$ws.Range("E424").Value = 11807
$ws.Range("F424").Value = 221340
$ws.Range("G424").Value = 443
